# Update main.py to provinces level:
# Add five new rows (Baoding city + its counties) below the existing
# Shijiazhuang rows, matching the look/format of the rows already there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 15) down onto the
# five new rows (16-20) so fonts/borders/alignment match the rest of the
# table, then overwrite the row-16 placeholder blank row in the process.
$ws.Range("A15:C15").Copy()
$ws.Range("A16:C20").PasteSpecial(-4122)
$ws.Range("A16:C20").RowHeight = 21.75

$data = @(
    @("保定市", "涞水县", "https://www.laishui.gov.cn/index.php?m=content&c=index&a=lists&catid=173"),
    @("保定市", "阜平县", "https://www.bdfuping.gov.cn/news/53/#c_news_list-1548139445064-1"),
    @("保定市", "定兴县", "http://www.dingxing.gov.cn/czyslist-394-more.html"),
    @("保定市", "唐县", "http://www.tangxian.gov.cn/czyslist-1116-more.html"),
    @("保定市", "高阳县", "https://www.gaoyang.gov.cn/cai/")
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
